$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is added at the top (row 2). Push the existing
# records down by one row first (copy row N into row N+1, working from the
# bottom up so no data is overwritten before it's been copied), which turns
# the old row 20 into the new row 21.
$ws.Range("A20:T20").Copy($ws.Range("A21:T21"))
$ws.Range("A19:T19").Copy($ws.Range("A20:T20"))
$ws.Range("A18:T18").Copy($ws.Range("A19:T19"))
$ws.Range("A17:T17").Copy($ws.Range("A18:T18"))
$ws.Range("A16:T16").Copy($ws.Range("A17:T17"))
$ws.Range("A15:T15").Copy($ws.Range("A16:T16"))
$ws.Range("A14:T14").Copy($ws.Range("A15:T15"))
$ws.Range("A13:T13").Copy($ws.Range("A14:T14"))
$ws.Range("A12:T12").Copy($ws.Range("A13:T13"))
$ws.Range("A11:T11").Copy($ws.Range("A12:T12"))
$ws.Range("A10:T10").Copy($ws.Range("A11:T11"))
$ws.Range("A9:T9").Copy($ws.Range("A10:T10"))
$ws.Range("A8:T8").Copy($ws.Range("A9:T9"))
$ws.Range("A7:T7").Copy($ws.Range("A8:T8"))
$ws.Range("A6:T6").Copy($ws.Range("A7:T7"))
$ws.Range("A5:T5").Copy($ws.Range("A6:T6"))
$ws.Range("A4:T4").Copy($ws.Range("A5:T5"))
$ws.Range("A3:T3").Copy($ws.Range("A4:T4"))
$ws.Range("A2:T2").Copy($ws.Range("A3:T3"))

# Now overwrite row 2 with the new record's data (same market/product/etc.,
# new date, volume and prices).
$ws.Range("D2").Value = 44956
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 3000
$ws.Range("O2").Value = 3000
$ws.Range("P2").Value = 3000
$ws.Range("Q2").Value = "$/bandeja 2 kilos"
$ws.Range("R2").Value = "Provincia de Diguillín"
$ws.Range("S2").Value = 1500
$ws.Range("T2").Value = 2
